{"js": "// Edit: turn \"...${nombre_solicitado}, de ${solicitado_tipo_edad} a\u00f1os de edad, ${solicitado_estado_civil}...\"\n// into \"...${nombre_solicitado},  ${solicitado_tipo_edad}, ${solicitado_estado_civil}...\"\n// i.e. drop the literal \"de \" run before the ${solicitado_tipo_edad} bookmark and the\n// literal \" a\u00f1os de edad\" run after it (both hard-coded words replaced by the\n// ${solicitado_tipo_edad} merge-field value itself), leaving a single space in their place.\n\nconst body = context.document.body;\n\n// 1) Locate the unique sentence fragment that anchors this edit.\nconst anchor = body.search(\"de ${solicitado_tipo_edad} a\u00f1os de edad\", { matchCase: true });\nanchor.load(\"items\");\nawait context.sync();\n\nif (anchor.items.length === 0) {\n  throw new Error(\"Could not find the '${solicitado_tipo_edad}' sentence to edit.\");\n}\nconst anchorRange = anchor.items[0];\n\n// 2) Remove the literal \"de \" that precedes the bookmark / merge field.\nconst prefix = anchorRange.search(\"de \", { matchCase: true });\nprefix.load(\"items\");\nawait context.sync();\nprefix.items[0].delete();\nawait context.sync();\n\n// 3) Remove the literal \" a\u00f1os de edad\" that follows the merge field.\nconst afterPrefixRemoved = body.search(\"${solicitado_tipo_edad} a\u00f1os de edad\", { matchCase: true });\nafterPrefixRemoved.load(\"items\");\nawait context.sync();\nconst suffix = afterPrefixRemoved.items[0].search(\" a\u00f1os de edad\", { matchCase: true });\nsuffix.load(\"items\");\nawait context.sync();\nsuffix.items[0].delete();\nawait context.sync();\n\n// 4) Re-insert a single space where \"de \" used to be, right after the _GoBack\n//    bookmark and before the ${solicitado_tipo_edad} field, so the sentence still\n//    reads \"... ${nombre_solicitado},  ${solicitado_tipo_edad}, ${solicitado_estado_civil} ...\".\nconst bookmark = context.document.getBookmarkRange(\"_GoBack\");\nbookmark.insertText(\" \", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Edit: turn \"...${nombre_solicitado}, de ${solicitado_tipo_edad} a\u00f1os de edad, ${solicitado_estado_civil}...\"\n# into \"...${nombre_solicitado},  ${solicitado_tipo_edad}, ${solicitado_estado_civil}...\"\n# i.e. drop the literal \"de \" run before the ${solicitado_tipo_edad} bookmark/merge-field\n# and the literal \" a\u00f1os de edad\" run after it, leaving a single space in their place.\n\n$d = $word.ActiveDocument\n\n# 1) Locate the unique sentence fragment that anchors this edit.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Text = \"de `${solicitado_tipo_edad} a\u00f1os de edad\"\n$found = $find.Execute()\nif (-not $found) {\n  throw \"Could not find the '`${solicitado_tipo_edad}' sentence to edit.\"\n}\n$anchorStart = $find.Parent.Start\n\n# 2) Remove the literal \"de \" that precedes the bookmark / merge field.\n$prefixRange = $d.Range($anchorStart, $anchorStart + 3)\n$prefixRange.Delete()\n\n# 3) Remove the literal \" a\u00f1os de edad\" that follows the merge field.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.MatchCase = $true\n$find2.Text = \"`${solicitado_tipo_edad}\"\n$find2.Execute() | Out-Null\n$varEnd = $find2.Parent.End\n\n$suffixText = \" a\u00f1os de edad\"\n$suffixRange = $d.Range($varEnd, $varEnd + $suffixText.Length)\n$suffixRange.Delete()\n\n# 4) Re-insert a single space where \"de \" used to be, right after the _GoBack\n#    bookmark and before the ${solicitado_tipo_edad} field, so the sentence still\n#    reads \"... ${nombre_solicitado},  ${solicitado_tipo_edad}, ${solicitado_estado_civil} ...\".\n$bookmark = $d.Bookmarks(\"_GoBack\")\n$insertPoint = $d.Range($bookmark.Start, $bookmark.Start)\n$insertPoint.InsertAfter(\" \")\n"}
